# Group_G_WDD_project_task_list - second checkpoint commit saturday session
# Applies the "task progress" and "lists" sheet updates described in the target diff.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "task progress"
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows above the old row 31 so everything from the old
# row 31 ("Let Emer know...") down to the end shifts down by two rows
# (old 31->33, old 33->35, old 34..41->36..43, old 44,45->46,47).
$ws.Range("A31:A32").EntireRow.Insert()

# Remove the old "Name / task assigned" mini-table that lived in columns
# F:G (including its data validation list), it is not part of the sheet
# any more.
$ws.Range("F1:G1").EntireColumn.Delete()

# --- In-place content updates (rows 1-21 keep their row numbers) ------
$ws.Range("D21").Value = "Alexander, Tristan"

# Row 22: "design board" -> "design board updated", mark done, assign tristan
$ws.Range("B22").Value = "design board updated"
$ws.Range("C22").Value = "y"
$ws.Range("D22").Value = "tristan"

# Row 23: navbar improvement now marked done or/assigned
$ws.Range("C23").Value = "y"
$ws.Range("D23").Value = "Alexander"

# Row 24 (previously empty) - new task
$ws.Range("B24").Value = "procedural navbar (change one file change it on every page)"
$ws.Range("C24").Value = "y"
$ws.Range("D24").Value = "Andrei"

# Row 26: append note about everyone present
$ws.Range("B26").Value = "Let Emer know if someone is out of contact"
$ws.Range("D26").Value = "everyone present and contributing as of 22/11/2025"

# Rows 29/30 (previously empty, now inside week 11 block)
$ws.Range("B29").Value = "fnctionality"
$ws.Range("D29").Value = "All"
$ws.Range("B30").Value = "form validation on each page"
$ws.Range("D30").Value = "All"

# Row 47 (old row 45, shifted +2): drop the stray "Alex" that used to sit
# in column G - already removed above along with the whole F:G columns.
$ws.Range("C47").Value = "Andrei"
$ws.Range("D47").Value = "Derek"
$ws.Range("E47").Value = "Tristan"

# ----------------------------------------------------------------------
# Sheet 2: "lists"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("I1").Value = $null
$ws2.Range("I2").Value = "WEEK 9"
$ws2.Range("K2").Value = "COMPLETED "
$ws2.Range("I3").Value = "WEEK 10"
$ws2.Range("K3").Value = "IN-progress"
$ws2.Range("I4").Value = "WEEK 11"
$ws2.Range("K4").Value = "to be started "
$ws2.Range("I5").Value = "WEEK 12"
$ws2.Range("K5").Value = "more time needed "
$ws2.Range("I6").Value = "WEEK 13"
$ws2.Range("I7").Value = "WEEK 14"
$ws2.Range("I8").Value = "WEEK 15"
$ws2.Range("I9").Value = "WEEK 16"

$ws2.Columns.Item(11).ColumnWidth = 12.85546875
